$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.235000252723694
$ws.Range("B1").Value = 1.90663206577301
$ws.Range("C1").Value = 2.455095529556274
$ws.Range("D1").Value = 3.857033491134644
$ws.Range("E1").Value = 1.154526948928833
